$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the name bug: "Hokie Bird" -> "Joe", and update his ticket count
$ws.Range("B2").Value = "Joe"
$ws.Range("D2").Value = 2

# Add the missing row for Bob
$ws.Range("A3").Value = 1234
$ws.Range("B3").Value = "Bob"
$ws.Range("D3").Value = 1

# Reset row heights back to the sheet's default
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

# Move the active selection
$ws.Range("D6").Select()
